$d = $word.ActiveDocument

# 1. Research Project Title
$d.Content.Find.Execute("Sample Research Title", $true, $true, $false, $false, $false, $true, 1, $false, "Web-Based Document Management System for Research Evaluation and Monitoring Center", 2)

# 2 & 3. Implementing Institution / Collaborating Institution (both occurrences of the same text)
$d.Content.Find.Execute("National Basketball Association", $true, $true, $false, $false, $false, $true, 1, $false, "Research Evaluation and Monitoring Center", 2)

# 4. Name of Project Leader (case-sensitive match so the all-caps signature line is untouched)
$d.Content.Find.Execute("John Doe", $true, $true, $false, $false, $false, $true, 1, $false, "Nicole Franzyne Jao", 2)

# 5. Project Duration
$d.Content.Find.Execute("1 year", $true, $true, $false, $false, $false, $true, 1, $false, "10 months", 2)

# 6. Team Members: rename the first two team members, then drop the third entirely
#    (original = "Lebron James" <br> "Anthony Davis" <br> "Austin Reaves")
$d.Content.Find.Execute("Lebron James", $true, $true, $false, $false, $false, $true, 1, $false, "Denice Shanley Alemania", 2)
$d.Content.Find.Execute("Anthony Davis", $true, $true, $false, $false, $false, $true, 1, $false, "Regienald Pueblos", 2)

$rng = $d.Content
$found = $rng.Find.Execute("Austin Reaves", $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) {
    $delRng = $d.Range($rng.Start - 1, $rng.End)
    $delRng.Delete()
}

# 7. Total Project Cost
$d.Content.Find.Execute("1,000,000", $true, $true, $false, $false, $false, $true, 1, $false, "0", 2)

# 8. Signature block name (uppercase form)
$d.Content.Find.Execute("JOHN DOE", $true, $true, $false, $false, $false, $true, 1, $false, "NICOLE FRANZYNE JAO", 2)
